$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial" (2nd Partial) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 2
$ws2.Range("E2").Value = 24
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 100
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 8.800000000000001
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0

# Row 3
$ws2.Range("E3").Value = 24
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 8.800000000000001
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

# Row 4
$ws2.Range("E4").Value = 37
$ws2.Range("F4").Value = 2
$ws2.Range("G4").Value = 94.90000000000001
$ws2.Range("H4").Value = 5.1
$ws2.Range("I4").Value = 8
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0

# Row 5
$ws2.Range("E5").Value = 37
$ws2.Range("F5").Value = 2
$ws2.Range("G5").Value = 94.90000000000001
$ws2.Range("H5").Value = 5.1
$ws2.Range("I5").Value = 8
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0

# Row 6
$ws2.Range("E6").Value = 24
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 100
$ws2.Range("H6").Value = 0
$ws2.Range("I6").Value = 9
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0

# Row 7
$ws2.Range("E7").Value = 24
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 100
$ws2.Range("H7").Value = 0
$ws2.Range("I7").Value = 9
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0

# Row 8
$ws2.Range("E8").Value = 22
$ws2.Range("F8").Value = 17
$ws2.Range("G8").Value = 56.4
$ws2.Range("H8").Value = 43.6
$ws2.Range("I8").Value = 9.9
$ws2.Range("J8").Value = 17
$ws2.Range("K8").Value = 43.59

# Row 9
$ws2.Range("E9").Value = 22
$ws2.Range("F9").Value = 17
$ws2.Range("G9").Value = 56.4
$ws2.Range("H9").Value = 43.6
$ws2.Range("I9").Value = 9.9
$ws2.Range("J9").Value = 17
$ws2.Range("K9").Value = 43.6

# Row 10
$ws2.Range("E10").Value = 107
$ws2.Range("F10").Value = 19
$ws2.Range("G10").Value = 84.90000000000001
$ws2.Range("H10").Value = 15.1
$ws2.Range("I10").Value = 8.9
$ws2.Range("J10").Value = 17
$ws2.Range("K10").Value = 13.5

# --- Sheet "Final" ---
$ws3 = $wb.Worksheets.Item("Final")

$ws3.Range("I2").Value = 8.800000000000001
$ws3.Range("I3").Value = 8.800000000000001
$ws3.Range("I4").Value = 8.199999999999999
$ws3.Range("I5").Value = 8.199999999999999
$ws3.Range("I6").Value = 9
$ws3.Range("I7").Value = 9
$ws3.Range("I8").Value = 9.1
$ws3.Range("I9").Value = 9.1
$ws3.Range("I10").Value = 8.800000000000001
